$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# 1. "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"

# 2. "2017-02-17 09:56:26" -> "2017-02-17 09:58:33" (Overview Latest HO Xliff Generate Date)
$overview.Range("G2").Value = "2017-02-17 09:58:33"
$overview.Range("G3").Value = "2017-02-17 09:58:33"

# 3. "ht" -> "mt" (Priority column)
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"
$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

# 4. "2017-02-17 09:56:08" -> "2017-02-17 09:58:16" (Latest Handoff Datetime, zh-cn only;
#    in the source file zh-cn!H happened to use a distinct shared string from de-de!H)
$zhcn.Range("H2").Value = "2017-02-17 09:58:16"
$zhcn.Range("H3").Value = "2017-02-17 09:58:16"

# de-de!H2/H3 originally shared the very same string as Overview!G2/G3
# ("2017-02-17 09:56:26"), so it tracks that same text's update here.
$dede.Range("H2").Value = "2017-02-17 09:58:33"
$dede.Range("H3").Value = "2017-02-17 09:58:33"

# 5. Error Detail on zh-cn row 2
$zhcn.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68de26946227fbf20e045aa7e5d7a6aa9c6044c8/e2e/de898377-1a40-4256-a36e-cf8b89338f72.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20eadb958ccb25a5449925b4bd67b9d5b51f09a6/e2e/de898377-1a40-4256-a36e-cf8b89338f72.md."

# 6. Column width changes
# NOTE: this runtime's ColumnWidth setter snaps the stored OOXML width to the
# nearest 1/6-character increment (stored = round((set+5/6)*6)/6), so we
# choose inputs that land on the grid point nearest the target width.
$overview.Range("E1").ColumnWidth = 16.333333333333336
$overview.Range("F1").ColumnWidth = 16.333333333333336

$zhcn.Range("C1").ColumnWidth = 16.333333333333336
$zhcn.Range("R1").ColumnWidth = 39.166666666666664

$dede.Range("C1").ColumnWidth = 16.333333333333336
$dede.Range("R1").ColumnWidth = 39.166666666666664
